# Applies the "Updated cryptos list" crypto-price refresh to Sheet1.
# Plain text updates use .Value directly. Updates whose new text reads as a
# number (e.g. "246.98") briefly force the cell to Text format ("@") so Excel
# stores the literal string instead of re-parsing it as a number, then restore
# the cell's original NumberFormat/Style so no formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '26.501.06'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.734.14'
$fmt = $ws.Range("D5").NumberFormat
$sty = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.98'
$ws.Range("D5").NumberFormat = $fmt
$ws.Range("D5").Style = $sty
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("E6").Value = '  +0.05%  '
$fmt = $ws.Range("D7").NumberFormat
$sty = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4891'
$ws.Range("D7").NumberFormat = $fmt
$ws.Range("D7").Style = $sty
$ws.Range("E7").Value = '  +1.70%  '
$fmt = $ws.Range("D8").NumberFormat
$sty = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2663'
$ws.Range("D8").NumberFormat = $fmt
$ws.Range("D8").Style = $sty
$ws.Range("E8").Value = '  -0.70%  '
$fmt = $ws.Range("D9").NumberFormat
$sty = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06219'
$ws.Range("D9").NumberFormat = $fmt
$ws.Range("D9").Style = $sty
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("D10").Value = '1.729.90'
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("E12").Value = '  -1.24%  '
$fmt = $ws.Range("D13").NumberFormat
$sty = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.589'
$ws.Range("D13").NumberFormat = $fmt
$ws.Range("D13").Style = $sty
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("E14").Value = '  -2.21%  '
$fmt = $ws.Range("D15").NumberFormat
$sty = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.41'
$ws.Range("D15").NumberFormat = $fmt
$ws.Range("D15").Style = $sty
$ws.Range("E15").Value = '  +0.35%  '
$fmt = $ws.Range("D16").NumberFormat
$sty = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").NumberFormat = $fmt
$ws.Range("D16").Style = $sty
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$fmt = $ws.Range("D17").NumberFormat
$sty = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007371'
$ws.Range("D17").NumberFormat = $fmt
$ws.Range("D17").Style = $sty
$ws.Range("E17").Value = '  +6.91%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '26.502.36'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("D21").Value = '1.953.18'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("E22").Value = '  -0.70%  '
$fmt = $ws.Range("D23").NumberFormat
$sty = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.746'
$ws.Range("D23").NumberFormat = $fmt
$ws.Range("D23").Style = $sty
$ws.Range("E23").Value = '  -1.83%  '
$fmt = $ws.Range("D25").NumberFormat
$sty = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.07'
$ws.Range("D25").NumberFormat = $fmt
$ws.Range("D25").Style = $sty
$ws.Range("E25").Value = '  +3.94%  '
$fmt = $ws.Range("D26").NumberFormat
$sty = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.43'
$ws.Range("D26").NumberFormat = $fmt
$ws.Range("D26").Style = $sty
$ws.Range("E26").Value = '  +0.27%  '
$fmt = $ws.Range("D27").NumberFormat
$sty = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.412'
$ws.Range("D27").NumberFormat = $fmt
$ws.Range("D27").Style = $sty
$ws.Range("E27").Value = '  -1.01%  '
$fmt = $ws.Range("D28").NumberFormat
$sty = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.771'
$ws.Range("D28").NumberFormat = $fmt
$ws.Range("D28").Style = $sty
$ws.Range("E28").Value = '  -2.29%  '
$fmt = $ws.Range("D29").NumberFormat
$sty = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.79'
$ws.Range("D29").NumberFormat = $fmt
$ws.Range("D29").Style = $sty
$ws.Range("E29").Value = '  +0.84%  '
$fmt = $ws.Range("D30").NumberFormat
$sty = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.015'
$ws.Range("D30").NumberFormat = $fmt
$ws.Range("D30").Style = $sty
$ws.Range("E30").Value = '  +0.06%  '
$fmt = $ws.Range("D31").NumberFormat
$sty = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08041'
$ws.Range("D31").NumberFormat = $fmt
$ws.Range("D31").Style = $sty
$ws.Range("E31").Value = '  +1.85%  '
$fmt = $ws.Range("D32").NumberFormat
$sty = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.687'
$ws.Range("D32").NumberFormat = $fmt
$ws.Range("D32").Style = $sty
$ws.Range("E32").Value = '  -1.52%  '
$fmt = $ws.Range("D33").NumberFormat
$sty = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04565'
$ws.Range("D33").NumberFormat = $fmt
$ws.Range("D33").Style = $sty
$ws.Range("E33").Value = '  -0.54%  '
$fmt = $ws.Range("D34").NumberFormat
$sty = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.000'
$ws.Range("D34").NumberFormat = $fmt
$ws.Range("D34").Style = $sty
$ws.Range("E34").Value = '  +0.07%  '
$fmt = $ws.Range("D35").NumberFormat
$sty = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.614'
$ws.Range("D35").NumberFormat = $fmt
$ws.Range("D35").Style = $sty
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +0.52%  '
$fmt = $ws.Range("D37").NumberFormat
$sty = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6339'
$ws.Range("D37").NumberFormat = $fmt
$ws.Range("D37").Style = $sty
$ws.Range("E37").Value = '  -0.54%  '
$fmt = $ws.Range("D38").NumberFormat
$sty = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8969'
$ws.Range("D38").NumberFormat = $fmt
$ws.Range("D38").Style = $sty
$ws.Range("E38").Value = '  -3.70%  '
$fmt = $ws.Range("D39").NumberFormat
$sty = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.020'
$ws.Range("D39").NumberFormat = $fmt
$ws.Range("D39").Style = $sty
$ws.Range("E39").Value = '  +1.16%  '
$fmt = $ws.Range("D40").NumberFormat
$sty = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.396'
$ws.Range("D40").NumberFormat = $fmt
$ws.Range("D40").Style = $sty
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  +0.17%  '
$fmt = $ws.Range("D42").NumberFormat
$sty = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01501'
$ws.Range("D42").NumberFormat = $fmt
$ws.Range("D42").Style = $sty
$ws.Range("E42").Value = '  -1.04%  '
$fmt = $ws.Range("D43").NumberFormat
$sty = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.75'
$ws.Range("D43").NumberFormat = $fmt
$ws.Range("D43").Style = $sty
$ws.Range("E43").Value = '  -9.03%  '
$fmt = $ws.Range("D44").NumberFormat
$sty = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.399'
$ws.Range("D44").NumberFormat = $fmt
$ws.Range("D44").Style = $sty
$ws.Range("E44").Value = '  -5.85%  '
$fmt = $ws.Range("D45").NumberFormat
$sty = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3885'
$ws.Range("D45").NumberFormat = $fmt
$ws.Range("D45").Style = $sty
$ws.Range("E45").Value = '  -0.92%  '
$fmt = $ws.Range("D46").NumberFormat
$sty = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.934'
$ws.Range("D46").NumberFormat = $fmt
$ws.Range("D46").Style = $sty
$ws.Range("E46").Value = '  -0.37%  '
$fmt = $ws.Range("D47").NumberFormat
$sty = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1184'
$ws.Range("D47").NumberFormat = $fmt
$ws.Range("D47").Style = $sty
$ws.Range("E47").Value = '  -1.19%  '
$fmt = $ws.Range("D48").NumberFormat
$sty = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05398'
$ws.Range("D48").NumberFormat = $fmt
$ws.Range("D48").Style = $sty
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("E49").Value = '  -1.17%  '
$fmt = $ws.Range("D50").NumberFormat
$sty = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.48'
$ws.Range("D50").NumberFormat = $fmt
$ws.Range("D50").Style = $sty
$ws.Range("E50").Value = '  -1.34%  '
$fmt = $ws.Range("D51").NumberFormat
$sty = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.256'
$ws.Range("D51").NumberFormat = $fmt
$ws.Range("D51").Style = $sty
$ws.Range("E51").Value = '  -0.24%  '
